$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells stay as text (avoid Excel auto-converting numeric-looking
# strings like "48.80" or "27.339.87" into numbers and losing formatting/precision).
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "B22", "C22", "D22", "E22", "B23", "C23", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.339.87"
$ws.Range("E2").Value = "  -3.78%  "
$ws.Range("D3").Value = "1.846.10"
$ws.Range("E3").Value = "  -5.66%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "321.24"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "0.4444"
$ws.Range("E7").Value = "  -6.65%  "
$ws.Range("D8").Value = "0.3819"
$ws.Range("E8").Value = "  -5.81%  "
$ws.Range("D9").Value = "48.80"
$ws.Range("E9").Value = "  -8.90%  "
$ws.Range("D10").Value = "0.07794"
$ws.Range("E10").Value = "  -8.03%  "
$ws.Range("D11").Value = "1.013"
$ws.Range("E11").Value = "  -4.23%  "
$ws.Range("D12").Value = "21.32"
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("D13").Value = "1.867.18"
$ws.Range("E13").Value = "  -4.51%  "
$ws.Range("D14").Value = "5.823"
$ws.Range("E14").Value = "  -5.11%  "
$ws.Range("D15").Value = "7.058"
$ws.Range("E15").Value = "  -7.27%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "85.08"
$ws.Range("E17").Value = "  -5.11%  "
$ws.Range("E18").Value = "  -4.49%  "
$ws.Range("D19").Value = "0.06505"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "16.89"
$ws.Range("E20").Value = "  -9.19%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.367.10"
$ws.Range("E22").Value = "  -3.71%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.439"
$ws.Range("E23").Value = "  -6.37%  "
$ws.Range("D24").Value = "10.73"
$ws.Range("E24").Value = "  -7.00%  "
$ws.Range("D25").Value = "2.259"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").Value = "2.062.40"
$ws.Range("E26").Value = "  -5.77%  "
$ws.Range("D27").Value = "151.53"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").Value = "19.24"
$ws.Range("E28").Value = "  -4.71%  "
$ws.Range("D29").Value = "2.040"
$ws.Range("E29").Value = "  -5.40%  "
$ws.Range("D30").Value = "5.466"
$ws.Range("E30").Value = "  -7.62%  "
$ws.Range("D31").Value = "119.67"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("D32").Value = "0.09304"
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("D33").Value = "1.459"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "0.9248"
$ws.Range("E34").Value = "  -5.41%  "
$ws.Range("D35").Value = "3.626"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "5.203"
$ws.Range("E36").Value = "  -6.87%  "
$ws.Range("D37").Value = "0.02209"
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("D38").Value = "0.05930"
$ws.Range("E38").Value = "  -4.36%  "
$ws.Range("D39").Value = "1.203"
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("D40").Value = "8.267"
$ws.Range("E40").Value = "  -6.50%  "
$ws.Range("D41").Value = "1.002"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "0.5903"
$ws.Range("E42").Value = "  -4.73%  "
$ws.Range("D43").Value = "0.1844"
$ws.Range("E43").Value = "  -3.92%  "
$ws.Range("D44").Value = "10.24"
$ws.Range("E44").Value = "  -7.90%  "
$ws.Range("D45").Value = "1.258"
$ws.Range("E45").Value = "  -5.92%  "
$ws.Range("D46").Value = "0.5633"
$ws.Range("E46").Value = "  -5.46%  "
$ws.Range("D47").Value = "12.11"
$ws.Range("E47").Value = "  -7.11%  "
$ws.Range("D48").Value = "3.350"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").Value = "1.909"
$ws.Range("E49").Value = "  -7.33%  "
$ws.Range("D50").Value = "0.06857"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").Value = "107.74"
$ws.Range("E51").Value = "  -2.38%  "
